$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D, E, G (rows 2-51) hold text values (prices, percentages, counters).
# Force Text format before assignment so Excel does not auto-coerce the
# numeric-looking strings into real numbers, then restore the default style
# so no stray formatting is introduced.
$affected = $ws.Range("D2:G51")
$affected.NumberFormat = "@"

$ws.Range("D2").Value = '310.38'
$ws.Range("E2").Value = '1.89%'
$ws.Range("G2").Value = '7'
$ws.Range("D3").Value = '38.98'
$ws.Range("E3").Value = '8.95%'
$ws.Range("G3").Value = '7'
$ws.Range("D4").Value = '5.113'
$ws.Range("E4").Value = '1.45%'
$ws.Range("G4").Value = '7'
$ws.Range("D5").Value = '0.08166'
$ws.Range("E5").Value = '2.58%'
$ws.Range("G5").Value = '7'
$ws.Range("D6").Value = '2.007'
$ws.Range("E6").Value = '4.79%'
$ws.Range("G6").Value = '7'
$ws.Range("D7").Value = '7.932'
$ws.Range("E7").Value = '1.95%'
$ws.Range("G7").Value = '7'
$ws.Range("D8").Value = '0.9341'
$ws.Range("G8").Value = '7'
$ws.Range("D9").Value = '0.1424'
$ws.Range("E9").Value = '8.65%'
$ws.Range("G9").Value = '7'
$ws.Range("D10").Value = '0.1959'
$ws.Range("E10").Value = '2.77%'
$ws.Range("G10").Value = '7'
$ws.Range("D11").Value = '0.09251'
$ws.Range("E11").Value = '1.38%'
$ws.Range("G11").Value = '7'
$ws.Range("D12").Value = '0.03472'
$ws.Range("E12").Value = '0.84%'
$ws.Range("G12").Value = '7'
$ws.Range("D13").Value = '0.09852'
$ws.Range("E13").Value = '0.08%'
$ws.Range("G13").Value = '7'
$ws.Range("D14").Value = '0.001407'
$ws.Range("E14").Value = '0.28%'
$ws.Range("G14").Value = '7'
$ws.Range("D15").Value = '0.005900'
$ws.Range("E15").Value = '-3.17%'
$ws.Range("G15").Value = '7'
$ws.Range("D16").Value = '3.571'
$ws.Range("E16").Value = '-4.17%'
$ws.Range("G16").Value = '7'
$ws.Range("D17").Value = '4.199'
$ws.Range("E17").Value = '1.82%'
$ws.Range("G17").Value = '7'
$ws.Range("D18").Value = '3.434'
$ws.Range("E18").Value = '1.15%'
$ws.Range("G18").Value = '7'
$ws.Range("D19").Value = '0.3449'
$ws.Range("E19").Value = '0.18%'
$ws.Range("G19").Value = '7'
$ws.Range("D20").Value = '0.1313'
$ws.Range("E20").Value = '0.26%'
$ws.Range("G20").Value = '7'
$ws.Range("D21").Value = '4.804'
$ws.Range("E21").Value = '-7.12%'
$ws.Range("G21").Value = '7'
$ws.Range("D22").Value = '0.2468'
$ws.Range("E22").Value = '4.99%'
$ws.Range("G22").Value = '7'
$ws.Range("D23").Value = '0.04453'
$ws.Range("E23").Value = '0.70%'
$ws.Range("G23").Value = '7'
$ws.Range("D24").Value = '0.001238'
$ws.Range("E24").Value = '0.36%'
$ws.Range("G24").Value = '7'
$ws.Range("E25").Value = '-9.82%'
$ws.Range("G25").Value = '7'
$ws.Range("G26").Value = '7'
$ws.Range("D27").Value = '0.0001303'
$ws.Range("E27").Value = '4.16%'
$ws.Range("G27").Value = '7'
$ws.Range("G28").Value = '7'
$ws.Range("G29").Value = '7'
$ws.Range("G30").Value = '7'
$ws.Range("G31").Value = '7'
$ws.Range("G32").Value = '7'
$ws.Range("G33").Value = '7'
$ws.Range("G34").Value = '7'
$ws.Range("G35").Value = '7'
$ws.Range("G36").Value = '7'
$ws.Range("G37").Value = '7'
$ws.Range("G38").Value = '7'
$ws.Range("D39").Value = '0.02131'
$ws.Range("E39").Value = '10.04%'
$ws.Range("G39").Value = '7'
$ws.Range("D40").Value = '0.05179'
$ws.Range("E40").Value = '-3.21%'
$ws.Range("G40").Value = '7'
$ws.Range("D41").Value = '0.007468'
$ws.Range("E41").Value = '-1.63%'
$ws.Range("G41").Value = '7'
$ws.Range("D42").Value = '0.009954'
$ws.Range("E42").Value = '-1.76%'
$ws.Range("G42").Value = '7'
$ws.Range("D43").Value = '0.1367'
$ws.Range("E43").Value = '1.09%'
$ws.Range("G43").Value = '7'
$ws.Range("D44").Value = '0.002134'
$ws.Range("E44").Value = '-0.78%'
$ws.Range("G44").Value = '7'
$ws.Range("D45").Value = '0.009745'
$ws.Range("E45").Value = '-3.70%'
$ws.Range("G45").Value = '7'
$ws.Range("D46").Value = '0.00006329'
$ws.Range("E46").Value = '2.39%'
$ws.Range("G46").Value = '7'
$ws.Range("E47").Value = '0.16%'
$ws.Range("G47").Value = '7'
$ws.Range("G48").Value = '7'
$ws.Range("D49").Value = '0.001603'
$ws.Range("E49").Value = '-3.37%'
$ws.Range("G49").Value = '7'
$ws.Range("D50").Value = '0.00002104'
$ws.Range("E50").Value = '0.16%'
$ws.Range("G50").Value = '7'
$ws.Range("D51").Value = '0.0002004'
$ws.Range("E51").Value = '0.16%'
$ws.Range("G51").Value = '7'

$affected.Style = "Normal"

